$d = $word.ActiveDocument

$replacements = @(
    @("2024-10-21 Monday", "2024-10-22 Tuesday"),
    @("805×7=5635", "855×6=5130"),
    @("799×5=3995", "545×6=3270"),
    @("290×2=580", "214×8=1712"),
    @("625×5=3125", "350×4=1400"),
    @("926×7=6482", "163×7=1141"),
    @("240×8=1920", "483×9=4347"),
    @("439×2=878", "361×7=2527"),
    @("841×2=1682", "274×6=1644"),
    @("464×6=2784", "625×7=4375"),
    @("719×8=5752", "309×3=927"),
    @("708×9=6372", "988×2=1976"),
    @("535×4=2140", "710×9=6390"),
    @("914×5=4570", "502×8=4016"),
    @("938×9=8442", "809×6=4854"),
    @("276×4=1104", "900×2=1800"),
    @("933×7=6531", "808×3=2424"),
    @("270×3=810", "892×7=6244"),
    @("533×4=2132", "780×4=3120"),
    @("181×9=1629", "310×5=1550"),
    @("629×6=3774", "272×5=1360"),
    @("122×6=732", "937×8=7496"),
    @("718×9=6462", "681×9=6129"),
    @("917×9=8253", "636×6=3816"),
    @("505×9=4545", "930×5=4650"),
    @("431×6=2586", "572×4=2288")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done replacing $($replacements.Count) text values"
